# Atualiza notas dos alunos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Edinaldo De Paiva Santos" entirely (row 4),
# which shifts the rows below it up by one.
$ws.Rows.Item(4).Delete()

# Fill in the grades (R1 column / column B) for the remaining students.
$ws.Range("B2").Value = 0   # Edgar Bispo Da Silva Neto
$ws.Range("B3").Value = 1   # Eduardo Lago Nunes
$ws.Range("B4").Value = 1   # Filipe Palma Abreu
$ws.Range("B5").Value = 0   # Jordan Santos Hohenfeld
$ws.Range("B6").Value = 1   # Lucas Borges Jagersbacher

# Update the active selection to match the target state.
$ws.Range("B7").Select()
